# Regenerate orders with updated distance/size codes.
# The experiment's distance conditions and one size code were renumbered:
#   D64 -> D69, D80 -> D86, D51 -> D55, S30 -> S31
# These codes are embedded as substrings throughout many shared strings
# (condition names, left/right filenames, and the standalone distance/size
# labels), so a substring Find & Replace across the whole used range
# updates every occurrence consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$rng = $ws.Cells

$rng.Replace("D64", "D69")
$rng.Replace("D80", "D86")
$rng.Replace("D51", "D55")
$rng.Replace("S30", "S31")
